# Insert two new price-report rows for "Vega Monumental Concepción - Cebolla"
# (weekly fruit/vegetable update). The new rows go in right above the
# existing "1a (guarda)" / "2a (guarda)" rows (old row 279), pushing all
# rows from 279 downward by two (old 279-292 become 281-294).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 279, shifting everything else (incl. the old
# 279-292 block) down to 281-294.
$ws.Rows("279:280").Insert()

# --- New row 279: "1a nueva(o)" ---
$ws.Range("A279").Value = 11
$ws.Range("B279").Value = 'Vega Monumental Concepción'
$ws.Range("C279").Value = 'Bíobío'
$ws.Range("D279").Value = 44509
$ws.Range("E279").Value = 8
$ws.Range("F279").Value = 100112004
$ws.Range("G279").Value = 'Cebolla'
$ws.Range("H279").Value = 'Sin especificar'
$ws.Range("I279").Value = '1a nueva(o)'
$ws.Range("J279").Value = 1000
$ws.Range("K279").Value = 2000
$ws.Range("L279").Value = 2200
$ws.Range("M279").Value = 2100
$ws.Range("N279").Value = '$/paquete 20 unidades (volumen en unidades)'
$ws.Range("O279").Value = "Región de O'Higgins"
$ws.Range("P279").Value = 105
$ws.Range("Q279").Value = 20
$ws.Range("R279").Value = 'Hortaliza'

# --- New row 280: "2a nueva(o)" ---
$ws.Range("A280").Value = 11
$ws.Range("B280").Value = 'Vega Monumental Concepción'
$ws.Range("C280").Value = 'Bíobío'
$ws.Range("D280").Value = 44509
$ws.Range("E280").Value = 8
$ws.Range("F280").Value = 100112004
$ws.Range("G280").Value = 'Cebolla'
$ws.Range("H280").Value = 'Sin especificar'
$ws.Range("I280").Value = '2a nueva(o)'
$ws.Range("J280").Value = 500
$ws.Range("K280").Value = 1800
$ws.Range("L280").Value = 1800
$ws.Range("M280").Value = 1800
$ws.Range("N280").Value = '$/paquete 20 unidades (volumen en unidades)'
$ws.Range("O280").Value = "Región de O'Higgins"
$ws.Range("P280").Value = 90
$ws.Range("Q280").Value = 20
$ws.Range("R280").Value = 'Hortaliza'
